# placeholder
